$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ON-PREM")
$ws.Range("D13").Value = 210000
